$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (capitalize / reword)
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Tiempo requerido"
$ws.Range("C1").Value = "Materias primas"
$ws.Range("D1").Value = "Precio de venta"

# Update the "materias primas" text for the last data row (pie de limon)
$ws.Range("C6").Value = "crema, limon,merengue,harina, huevos"

# Column widths (target stored widths: B=15.88671875, C=36.44140625, D=14;
# inputs below are tuned so the engine's stored `width` lands on the closest
# achievable value given its internal character-width quantization)
$ws.Columns.Item(2).ColumnWidth = 15.0
$ws.Columns.Item(3).ColumnWidth = 35.65
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666

# Selection
$ws.Range("C6").Select() | Out-Null

# Page setup
$ws.PageSetup.Orientation = 1
